$d = $word.ActiveDocument

$pairs = @(
    @("43×29=", "27×38="),
    @("94×67=", "97×33="),
    @("60×57=", "72×22="),
    @("15×29=", "34×82="),
    @("95×65=", "76×48="),
    @("62×68=", "90×65="),
    @("61×73=", "85×64="),
    @("50×52=", "66×23="),
    @("70×78=", "91×25="),
    @("31×58=", "49×39="),
    @("69×52=", "84×36="),
    @("18×25=", "27×46="),
    @("77×76=", "44×87="),
    @("27×40=", "60×71="),
    @("29×14=", "57×59="),
    @("20×81=", "33×90="),
    @("41×89=", "96×84="),
    @("34×23=", "46×29="),
    @("45×41=", "49×98="),
    @("34×30=", "50×28="),
    @("26×32=", "29×18="),
    @("88×76=", "80×94="),
    @("11×46=", "61×67="),
    @("73×59=", "15×41="),
    @("83×18=", "15×40=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
